$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: given a Range, forcibly materialise it as its own <w:r> (splitting
# it away from any neighbouring run it currently shares formatting with) by
# round-tripping its text through a distinct placeholder before writing the
# final text. This mirrors what Word does when a user selects a sub-string
# of a run and retypes it.
# ---------------------------------------------------------------------------
function Split-AndSetText($rng, [string]$finalText) {
    $placeholder = "@@@SPLIT@@@"
    $rng.Text = $placeholder
    $len = $finalText.Length
    $newEnd = $rng.Start + $len
    $rng2 = $d.Range($rng.Start, $rng.Start + $placeholder.Length)
    $rng2.Text = $finalText
    return $d.Range($rng.Start, $newEnd)
}

# ---------------------------------------------------------------------------
# 1) "GravimetricMoistureContent" -> "Gravimetric" + "Water" + "Content"
#    (3 occurrences; the middle one is bold+italic, the other two italic-only)
# ---------------------------------------------------------------------------
$searchStart = 0
while ($true) {
    $rng = $d.Range($searchStart, $d.Content.End)
    $found = $rng.Find.Execute("GravimetricMoistureContent")
    if (-not $found) { break }
    $matchStart = $rng.Start
    $matchEnd = $rng.End

    # "Moisture" sits between "Gravimetric" (11 chars) and "Content" (7 chars)
    $moistureRange = $d.Range($matchStart + 11, $matchStart + 19)
    $newRange = Split-AndSetText $moistureRange "Water"

    $searchStart = $newRange.End
}

# ---------------------------------------------------------------------------
# 2) "#VALUE!" -> "#" + "VALUE" (highlighted yellow) + "!"
# ---------------------------------------------------------------------------
$searchStart = 0
while ($true) {
    $rng = $d.Range($searchStart, $d.Content.End)
    $found = $rng.Find.Execute("#VALUE!")
    if (-not $found) { break }
    $matchStart = $rng.Start
    $matchEnd = $rng.End

    $valueRange = $d.Range($matchStart + 1, $matchStart + 6)
    $newValueRange = Split-AndSetText $valueRange "VALUE"
    $newValueRange.Font.HighlightColorIndex = 7

    $bangRange = $d.Range($newValueRange.End, $newValueRange.End + 1)
    $newBangRange = Split-AndSetText $bangRange "!"

    $searchStart = $newBangRange.End
}

# ---------------------------------------------------------------------------
# 3) "Gravimetric moisture content" -> "Gravimetric " + "Water" + " content"
#    (5 occurrences, all bold+italic)
# ---------------------------------------------------------------------------
$searchStart = 0
while ($true) {
    $rng = $d.Range($searchStart, $d.Content.End)
    $found = $rng.Find.Execute("Gravimetric moisture content")
    if (-not $found) { break }
    $matchStart = $rng.Start
    $matchEnd = $rng.End

    # "moisture" sits between "Gravimetric " (12 chars) and " content" (8 chars)
    $moistureRange = $d.Range($matchStart + 12, $matchStart + 20)
    $newRange = Split-AndSetText $moistureRange "Water"

    $searchStart = $newRange.End
}
